$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the existing row 34 (actc row), shifting the
# remaining rows (old 34-36) down to 35-37. This makes room for a new
# "sched_se_tax" / "Self-Employment Tax" entry.
$ws.Rows("34:34").Insert()

# The freshly inserted row 34 doesn't inherit the surrounding table's
# formatting automatically, so copy the cell formats from row 35 (the
# row that used to be row 34, directly below the new blank row) into the
# corresponding cells of row 34. Copy column-by-column so we don't touch
# column C (which has no content/style on this row in the target).
$ws.Range("A35").Copy()
$ws.Range("A34").PasteSpecial(-4122)
$ws.Range("B35").Copy()
$ws.Range("B34").PasteSpecial(-4122)
$ws.Range("D35").Copy()
$ws.Range("D34").PasteSpecial(-4122)
$ws.Range("E35").Copy()
$ws.Range("E34").PasteSpecial(-4122)

# Populate the new row with the variable name and readable name for the
# Self-Employment Tax line item.
$ws.Range("B34").Value = "sched_se_tax"
$ws.Range("D34").Value = "Self-Employment Tax"

# Match the author's final selection, which ended up resting on the new
# row's variable-name cell.
$ws.Range("B34").Select() | Out-Null
